$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.196.96"
$ws.Range("E2").Value = "  +5.32%  "
$ws.Range("D3").Value = "'2.741.46"
$ws.Range("E3").Value = "  +2.74%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'582.24"
$ws.Range("E5").Value = "  +2.75%  "
$ws.Range("D6").Value = "'158.61"
$ws.Range("E6").Value = "  +9.79%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +2.11%  "
$ws.Range("D9").Value = "'2.770.20"
$ws.Range("E9").Value = "  +3.80%  "
$ws.Range("E10").Value = "  +2.70%  "
$ws.Range("E11").Value = "  +6.27%  "
$ws.Range("E12").Value = "  +3.59%  "
$ws.Range("E13").Value = "  +2.04%  "
$ws.Range("D14").Value = "'3.232.68"
$ws.Range("E14").Value = "  +3.03%  "
$ws.Range("E15").Value = "  +4.29%  "
$ws.Range("D16").Value = "'64.117.29"
$ws.Range("E16").Value = "  +5.19%  "
$ws.Range("E17").Value = "  +8.18%  "
$ws.Range("D18").Value = "'2.764.54"
$ws.Range("E18").Value = "  +3.73%  "
$ws.Range("D19").Value = "'12.10"
$ws.Range("E19").Value = "  +4.51%  "
$ws.Range("D20").Value = "'4.96"
$ws.Range("E20").Value = "  +4.58%  "
$ws.Range("D21").Value = "'363.58"
$ws.Range("E21").Value = "  +3.50%  "
$ws.Range("D22").Value = "'7.06"
$ws.Range("E22").Value = "  +2.60%  "
$ws.Range("D23").Value = "'0.997"
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D24").Value = "'0.535"
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("D25").Value = "'67.41"
$ws.Range("E25").Value = "  +5.34%  "
$ws.Range("E26").Value = "  +5.89%  "
$ws.Range("D27").Value = "'8.63"
$ws.Range("E27").Value = "  +4.88%  "
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("D29").Value = "'0.0₃0922"
$ws.Range("E29").Value = "  +13.74%  "
$ws.Range("E30").Value = "  +1.25%  "
$ws.Range("D31").Value = "'7.22"
$ws.Range("E31").Value = "  +5.54%  "
$ws.Range("D32").Value = "'1.26"
$ws.Range("E32").Value = "  +19.13%  "
$ws.Range("E33").Value = "  +6.95%  "
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "'0.997"
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "'20.72"
$ws.Range("E35").Value = "  +3.93%  "
$ws.Range("E36").Value = "  +5.98%  "
$ws.Range("E38").Value = "  +9.96%  "
$ws.Range("E39").Value = "  +11.44%  "
$ws.Range("D40").Value = "'4.32"
$ws.Range("E40").Value = "  +5.33%  "
$ws.Range("D41").Value = "'341.42"
$ws.Range("E41").Value = "  +0.37%  "
$ws.Range("D42").Value = "'39.37"
$ws.Range("E42").Value = "  +2.20%  "
$ws.Range("D43").Value = "'5.87"
$ws.Range("E43").Value = "  +12.47%  "
$ws.Range("D44").Value = "'22.10"
$ws.Range("E44").Value = "  +8.48%  "
$ws.Range("D45").Value = "'22.19"
$ws.Range("E45").Value = "  +7.66%  "
$ws.Range("E46").Value = "  +6.63%  "
$ws.Range("D47").Value = "'0.654"
$ws.Range("E47").Value = "  +4.64%  "
$ws.Range("D48").Value = "'0.0261"
$ws.Range("E48").Value = "  +4.80%  "
$ws.Range("D49").Value = "'138.61"
$ws.Range("E49").Value = "  +4.17%  "
$ws.Range("E50").Value = "  +2.43%  "
$ws.Range("E51").Value = "  -0.04%  "
